# Scheduled runner refresh: update cached market-price / profit figures
# across the per-job Leve profit tables (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 221.125
$ws.Range("I12").Value = 229.0
$ws.Range("K12").Value = 229.0
$ws.Range("M12").Value = -59.0
$ws.Range("H17").Value = 7201.0
$ws.Range("I17").Value = 40000.0
$ws.Range("J17").Value = 4678.0
$ws.Range("K17").Value = 120000.0
$ws.Range("L17").Value = 14034.0
$ws.Range("M17").Value = -119832.0
$ws.Range("N17").Value = -14370.0
$ws.Range("H33").Value = 50170.4
$ws.Range("I33").Value = 50170.4
$ws.Range("J33").Value = 0.0
$ws.Range("K33").Value = 50170.4
$ws.Range("L33").Value = 0.0
$ws.Range("M33").Value = -49941.4
$ws.Range("H38").Value = 88.9
$ws.Range("I38").Value = 88.9
$ws.Range("K38").Value = 266.7
$ws.Range("M38").Value = 105.3
$ws.Range("H96").Value = 549.0769
$ws.Range("I96").Value = 382.66666
$ws.Range("J96").Value = 923.5
$ws.Range("K96").Value = 1147.99998
$ws.Range("L96").Value = 2770.5
$ws.Range("M96").Value = 225.0000199999999
$ws.Range("N96").Value = -5516.5
$ws.Range("H100").Value = 7404.4546
$ws.Range("I100").Value = 9408.167
$ws.Range("K100").Value = 9408.167
$ws.Range("M100").Value = -8867.167
$ws.Range("H103").Value = 460.66666
$ws.Range("I103").Value = 288.5
$ws.Range("J103").Value = 805.0
$ws.Range("K103").Value = 865.5
$ws.Range("L103").Value = 2415.0
$ws.Range("M103").Value = -279.5
$ws.Range("N103").Value = -3587.0
$ws.Range("H116").Value = 48717156.0
$ws.Range("I116").Value = 50221640.0
$ws.Range("K116").Value = 50221640.0
$ws.Range("M116").Value = -50218198.0
$ws.Range("H132").Value = 3662.3809
$ws.Range("I132").Value = 3327.2632
$ws.Range("K132").Value = 9981.7896
$ws.Range("M132").Value = -7451.7896
$ws.Range("H135").Value = 38462896.0
$ws.Range("I135").Value = 45455810.0
$ws.Range("J135").Value = 1875.0
$ws.Range("K135").Value = 409102290.0
$ws.Range("L135").Value = 16875.0
$ws.Range("M135").Value = -409099755.0
$ws.Range("N135").Value = -21945.0

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2549.0466
$ws.Range("I32").Value = 2707.3076
$ws.Range("K32").Value = 2707.3076
$ws.Range("M32").Value = -2420.3076
$ws.Range("H45").Value = 2124.5
$ws.Range("I45").Value = 2027.2222
$ws.Range("K45").Value = 2027.2222
$ws.Range("M45").Value = -1650.2222
$ws.Range("H61").Value = 15876490.0
$ws.Range("I61").Value = 20836518.0
$ws.Range("K61").Value = 20836518.0
$ws.Range("M61").Value = -20836306.0
$ws.Range("H132").Value = 27029016.0
$ws.Range("I132").Value = 33335314.0
$ws.Range("J132").Value = 2025.2858
$ws.Range("K132").Value = 100005942.0
$ws.Range("L132").Value = 6075.857400000001
$ws.Range("M132").Value = -100003412.0
$ws.Range("N132").Value = -11135.8574
$ws.Range("H136").Value = 15876490.0
$ws.Range("I136").Value = 20836518.0
$ws.Range("K136").Value = 62509554.0
$ws.Range("M136").Value = -62507004.0

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H7").Value = 2153.5
$ws.Range("J7").Value = 4004.0
$ws.Range("L7").Value = 4004.0
$ws.Range("N7").Value = -4230.0
$ws.Range("H87").Value = 0.0
$ws.Range("J87").Value = 0.0
$ws.Range("L87").Value = 0.0
$ws.Range("H90").Value = 0.0
$ws.Range("J90").Value = 0.0
$ws.Range("L90").Value = 0.0
$ws.Range("H105").Value = 1582.963
$ws.Range("I105").Value = 1492.1305
$ws.Range("K105").Value = 1492.1305
$ws.Range("M105").Value = 254.8695
$ws.Range("H107").Value = 71486856.0
$ws.Range("I107").Value = 41599.8
$ws.Range("K107").Value = 41599.8
$ws.Range("M107").Value = -39679.8
$ws.Range("H134").Value = 1217.3334
$ws.Range("I134").Value = 1161.4286
$ws.Range("J134").Value = 2000.0
$ws.Range("K134").Value = 3484.2858
$ws.Range("L134").Value = 6000.0
$ws.Range("M134").Value = -949.2857999999997
$ws.Range("N134").Value = -11070.0
$ws.Range("H138").Value = 0.0
$ws.Range("J138").Value = 0.0
$ws.Range("L138").Value = 0.0
$ws.Range("H139").Value = 44000.0
$ws.Range("I139").Value = 44000.0
$ws.Range("J139").Value = 0.0
$ws.Range("K139").Value = 44000.0
$ws.Range("L139").Value = 0.0
$ws.Range("M139").Value = -38860.0

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H51").Value = 74666.336
$ws.Range("I51").Value = 71999.75
$ws.Range("J51").Value = 79999.5
$ws.Range("K51").Value = 71999.75
$ws.Range("L51").Value = 79999.5
$ws.Range("M51").Value = -71263.75
$ws.Range("N51").Value = -81471.5
$ws.Range("H61").Value = 74666.336
$ws.Range("I61").Value = 71999.75
$ws.Range("J61").Value = 79999.5
$ws.Range("K61").Value = 71999.75
$ws.Range("L61").Value = 79999.5
$ws.Range("M61").Value = -71651.75
$ws.Range("N61").Value = -80695.5
$ws.Range("H105").Value = 2945.0
$ws.Range("I105").Value = 1477.1428
$ws.Range("J105").Value = 5000.0
$ws.Range("K105").Value = 1477.1428
$ws.Range("L105").Value = 5000.0
$ws.Range("M105").Value = 269.8571999999999
$ws.Range("N105").Value = -8494.0
$ws.Range("H107").Value = 1313.8462
$ws.Range("I107").Value = 1578.6
$ws.Range("J107").Value = 431.33334
$ws.Range("K107").Value = 1578.6
$ws.Range("L107").Value = 431.33334
$ws.Range("M107").Value = 341.4000000000001
$ws.Range("N107").Value = -4271.33334
$ws.Range("H122").Value = 2556.0667
$ws.Range("I122").Value = 2861.75
$ws.Range("K122").Value = 8585.25
$ws.Range("M122").Value = -6135.25
$ws.Range("H132").Value = 2171.1155
$ws.Range("I132").Value = 2060.4167
$ws.Range("J132").Value = 3499.5
$ws.Range("K132").Value = 6181.250100000001
$ws.Range("L132").Value = 10498.5
$ws.Range("M132").Value = -3651.250100000001
$ws.Range("N132").Value = -15558.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H69").Value = 4000.0
$ws.Range("J69").Value = 4000.0
$ws.Range("L69").Value = 12000.0
$ws.Range("N69").Value = -13622.0
$ws.Range("H72").Value = 4000.0
$ws.Range("J72").Value = 4000.0
$ws.Range("L72").Value = 36000.0
$ws.Range("N72").Value = -44112.0
$ws.Range("H98").Value = 800.0
$ws.Range("I98").Value = 0.0
$ws.Range("J98").Value = 800.0
$ws.Range("K98").Value = 0.0
$ws.Range("L98").Value = 2400.0
$ws.Range("N98").Value = -5396.0

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 650.0
$ws.Range("I97").Value = 650.0
$ws.Range("K97").Value = 650.0
$ws.Range("M97").Value = -154.0

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 71574456.0
$ws.Range("I7").Value = 83502696.0
$ws.Range("K7").Value = 83502696.0
$ws.Range("M7").Value = -83502584.0
$ws.Range("H86").Value = 0.0
$ws.Range("J86").Value = 0.0
$ws.Range("L86").Value = 0.0
$ws.Range("H89").Value = 0.0
$ws.Range("J89").Value = 0.0
$ws.Range("L89").Value = 0.0
$ws.Range("H93").Value = 1077.1333
$ws.Range("I93").Value = 1152.0834
$ws.Range("K93").Value = 1152.0834
$ws.Range("M93").Value = 95.91660000000002
$ws.Range("H100").Value = 3058.0588
$ws.Range("I100").Value = 2881.6365
$ws.Range("K100").Value = 2881.6365
$ws.Range("M100").Value = -2340.6365
$ws.Range("H126").Value = 71574456.0
$ws.Range("I126").Value = 83502696.0
$ws.Range("K126").Value = 250508088.0
$ws.Range("M126").Value = -250505618.0
$ws.Range("H132").Value = 6208.3477
$ws.Range("I132").Value = 3378.1428
$ws.Range("J132").Value = 10610.889
$ws.Range("K132").Value = 10134.4284
$ws.Range("L132").Value = 31832.667
$ws.Range("M132").Value = -7604.428400000001
$ws.Range("N132").Value = -36892.667

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 826.7143
$ws.Range("I107").Value = 826.7143
$ws.Range("K107").Value = 2480.1429
$ws.Range("M107").Value = -560.1428999999998
$ws.Range("H132").Value = 5178.4443
$ws.Range("I132").Value = 4944.476
$ws.Range("K132").Value = 14833.428
$ws.Range("M132").Value = -12303.428
$ws.Range("H133").Value = 126902.0
$ws.Range("J133").Value = 126902.0
$ws.Range("L133").Value = 126902.0
$ws.Range("N133").Value = -137022.0
$ws.Range("H136").Value = 2946.5386
$ws.Range("J136").Value = 8998.667
$ws.Range("L136").Value = 26996.001
$ws.Range("N136").Value = -32096.001

# Cells dropped from the refreshed export (no cached value for this run)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("N33").ClearContents()
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("N87").ClearContents()
$ws.Range("N90").ClearContents()
$ws.Range("N138").ClearContents()
$ws.Range("N139").ClearContents()
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("M98").ClearContents()
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("N86").ClearContents()
$ws.Range("N89").ClearContents()
